$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(" Dubai (DSC)", " October 27 2020", "Sunrisers won by 88 runs", "Sunrisers Hyderabad", "Delhi Capitals", "Wriddhiman Saha †", "87", "45", "12", "2", "193.33"),
    @(" Sharjah", " November 03 2020", "Sunrisers won by 10 wickets (with 17 balls remaining)", "Sunrisers Hyderabad", "Mumbai Indians", "Wriddhiman Saha †", "58", "45", "7", "1", "128.88"),
    @(" Abu Dhabi", " September 26 2020", "KKR won by 7 wickets (with 12 balls remaining)", "Sunrisers Hyderabad", "Kolkata Knight Riders", "Wriddhiman Saha †", "30", "31", "1", "1", "96.77"),
    @(" Sharjah", " October 31 2020", "Sunrisers won by 5 wickets (with 35 balls remaining)", "Sunrisers Hyderabad", "Royal Challengers Bangalore", "Wriddhiman Saha †", "39", "32", "4", "1", "121.87")
)

$startRow = 6
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    for ($c = 0; $c -lt $rowData.Count; $c++) {
        $col = $c + 1
        $cell = $ws.Cells.Item($row, $col)
        $text = $rowData[$c]
        # Columns G:K hold numeric-looking text (runs/balls/4s/6s/sr) that
        # must stay text, exactly like the rest of the sheet (t="str"
        # cells). Excel auto-converts a numeric-looking value to a real
        # number unless the cell is pre-formatted as Text, so force that
        # format before writing those values.
        if ($col -ge 7) {
            $cell.NumberFormat = "@"
        }
        $cell.Value = $text
    }
}
